$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 297-298 (everything from old row 297 onward
# shifts down by 2, old row 297 -> new row 299, ..., old row 316 -> new row 318)
$ws.Range("A297:A298").EntireRow.Insert()

# New row 297: Femacal de La Calera / Coquimbo / Ciruela / Angeleno / Primera
$ws.Range("A297").Value = 3
$ws.Range("B297").Value = "Femacal de La Calera"
$ws.Range("C297").Value = "Coquimbo"
$ws.Range("D297").Value = 45008
$ws.Range("E297").Value = 5
$ws.Range("F297").Value = "Fruta"
$ws.Range("G297").Value = 100103
$ws.Range("H297").Value = "Frutos de hueso (carozo)"
$ws.Range("I297").Value = 100103002
$ws.Range("J297").Value = "Ciruela"
$ws.Range("K297").Value = "Angeleno"
$ws.Range("L297").Value = "Primera"
$ws.Range("M297").Value = 56
$ws.Range("N297").Value = 10000
$ws.Range("O297").Value = 10000
$ws.Range("P297").Value = 10000
$ws.Range("Q297").Value = "$/caja 10 kilos"
$ws.Range("R297").Value = "Región de O'Higgins"
$ws.Range("S297").Value = 1000
$ws.Range("T297").Value = 10

# New row 298: Femacal de La Calera / Coquimbo / Ciruela / Angeleno / Primera
$ws.Range("A298").Value = 3
$ws.Range("B298").Value = "Femacal de La Calera"
$ws.Range("C298").Value = "Coquimbo"
$ws.Range("D298").Value = 45008
$ws.Range("E298").Value = 5
$ws.Range("F298").Value = "Fruta"
$ws.Range("G298").Value = 100103
$ws.Range("H298").Value = "Frutos de hueso (carozo)"
$ws.Range("I298").Value = 100103002
$ws.Range("J298").Value = "Ciruela"
$ws.Range("K298").Value = "Angeleno"
$ws.Range("L298").Value = "Primera"
$ws.Range("M298").Value = 40
$ws.Range("N298").Value = 8000
$ws.Range("O298").Value = 8000
$ws.Range("P298").Value = 8000
$ws.Range("Q298").Value = "$/caja 10 kilos"
$ws.Range("R298").Value = "Región de O'Higgins"
$ws.Range("S298").Value = 800
$ws.Range("T298").Value = 10
